# Generate Report for Handoff
#
# The localization-status report is regenerated: the "Status"/"zh-cn"/"de-de"
# columns move from "Handed back: in sync with en-US" to "Ready for handoff",
# and the handoff timestamps on the Overview sheet and the zh-cn sheet are
# refreshed to the new generation time. The Status column (now holding the
# shorter "Ready for handoff" text) is narrower, so its column is shrunk to
# match the new content.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns + HO Xliff generate date ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_overview.Range("G2").Value = "2016-08-27 00:59:15"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_zhcn.Range("H2").Value = "2016-08-27 00:59:10"

# --- de-de sheet: Status + Latest Handback DateTime ---
$ws_dede.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("H2").Value = "2016-08-27 00:59:15"

# --- Shrink the status columns now that they hold the shorter text ---
$ws_overview.Columns.Item(5).ColumnWidth = 16.29
$ws_overview.Columns.Item(6).ColumnWidth = 16.29
$ws_zhcn.Columns.Item(3).ColumnWidth = 16.29
$ws_dede.Columns.Item(3).ColumnWidth = 16.29
